# Update "Training Dashboard" sheet: new progress as of 04-Nov-2025.
# For each of rows 3-8, PERIOD TO EXPIRE (col H) decreases by 1 day and
# LAST UPDATE (col I) moves from 03-Nov-2025 to 04-Nov-2025.
#
# Column I holds the date as literal text (not a real date serial), so a
# plain `.Value = "04-Nov-2025"` assignment would get auto-recognized by
# Excel as a date and reformatted (changing the cell's style/number
# format). To keep the cell's existing style untouched we instead write
# the text via a formula (`="04-Nov-2025"`) — which Excel does not
# date-sniff — and then immediately convert that formula down to a
# plain value with copy / paste-special-values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-LiteralText($range, [string]$text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163) # xlPasteValues
}

$updates = @(
    @{ Row = 3; Period = 345; },
    @{ Row = 4; Period = 364; },
    @{ Row = 5; Period = 364; },
    @{ Row = 6; Period = -35; },
    @{ Row = 7; Period = 604; },
    @{ Row = 8; Period = 377; }
)

foreach ($u in $updates) {
    $row = $u.Row
    $ws.Range("H$row").Value = $u.Period
    Set-LiteralText $ws.Range("I$row") "04-Nov-2025"
}

$excel.CutCopyMode = $false
